# Updated cryptos list on Sun Oct 15 08:47:33 UTC 2023 with GitHub Actions
#
# Applies the per-row Price (D) and Volume(1h) (E) text updates described
# by the commit diff. Price values are written as exact text strings
# (forcing a "@" text number-format while assigning, then reverting the
# cell style to "Normal" so no stray style index is left behind) because
# several of them (e.g. "1.00", "22.09", "0.250", "7.40") would otherwise
# be auto-converted by Excel into numeric values and lose their original
# textual representation (trailing zeros, thousands separators, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# row -> (Price, Volume)
$rows = @{
    2  = @("27.054.73",  "  +0.56%  ")
    3  = @("1.565.35",   "  +1.03%  ")
    4  = @("1.00",       "  +0.32%  ")
    5  = @("208.82",     "  +1.12%  ")
    6  = @("0.491",      "  +0.88%  ")
    8  = @("22.09",      "  -0.28%  ")
    9  = @("0.250",      "  +1.37%  ")
    10 = @($null,        "  +1.97%  ")
    11 = @($null,        "  +0.43%  ")
    12 = @("1.787.69",   "  +0.97%  ")
    13 = @("1.556.62",   "  +0.44%  ")
    14 = @($null,        "  +0.60%  ")
    15 = @($null,        "  +0.67%  ")
    16 = @("27.035.53",  "  +0.49%  ")
    17 = @("61.93",      "  +0.52%  ")
    18 = @("0.0₃0706",   "  +1.18%  ")
    19 = @("216.36",     "  -0.45%  ")
    20 = @("7.40",       "  +2.01%  ")
    21 = @($null,        "  +0.34%  ")
    22 = @($null,        "  +2.40%  ")
    23 = @($null,        "  +0.12%  ")
    24 = @("1.95",       "  +0.03%  ")
    25 = @("154.09",     "  -0.10%  ")
    26 = @("6.62",       "  +0.03%  ")
    27 = @($null,        "  +1.14%  ")
    28 = @($null,        "  +1.85%  ")
    29 = @($null,        "  +0.29%  ")
    30 = @($null,        "  +1.64%  ")
    32 = @($null,        "  +0.25%  ")
    33 = @($null,        "  +4.83%  ")
    34 = @("1.427.05",   "  +0.84%  ")
    35 = @("1.09",       "  +12.58%  ")
    36 = @($null,        "  +2.01%  ")
    37 = @($null,        "  +2.70%  ")
    38 = @($null,        "  +1.48%  ")
    39 = @("0.533",      "  +1.59%  ")
    40 = @("0.813",      "  +0.79%  ")
    41 = @("5.80",       "  +1.71%  ")
    42 = @($null,        "  +0.39%  ")
    43 = @("2.32",       "  -0.19%  ")
    44 = @($null,        "  +0.15%  ")
    45 = @("64.82",      "  +0.57%  ")
    46 = @($null,        "  -0.07%  ")
    47 = @("1.700.47",   "  +0.96%  ")
    48 = @("86.64",      "  -1.02%  ")
    49 = @($null,        "  +1.79%  ")
    50 = @($null,        "  +0.25%  ")
    51 = @("0.0962",     "  +0.43%  ")
}

foreach ($row in ($rows.Keys | Sort-Object)) {
    $price  = $rows[$row][0]
    $volume = $rows[$row][1]

    if ($null -ne $price) {
        Set-TextValue "D$row" $price
    }
    if ($null -ne $volume) {
        $ws.Range("E$row").Value = $volume
    }
}
